$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty Embedding Time (D), LLM Response Time (E),
# Embedding Memory (G) and Query Memory (H) columns for rows 14-19.
# Columns F (Total Time) and I (Total Memory) already contain
# shared SUM formulas that will recalculate automatically.

$ws.Range("D14").Value = 26.73
$ws.Range("E14").Value = 9.73
$ws.Range("G14").Value = -16
$ws.Range("H14").Value = -47.44

$ws.Range("D15").Value = 23.34
$ws.Range("E15").Value = 40.71
$ws.Range("G15").Value = -17.7
$ws.Range("H15").Value = -47.55

$ws.Range("D16").Value = 36.450000000000003
$ws.Range("E16").Value = 36.479999999999997
$ws.Range("G16").Value = -14.67
$ws.Range("H16").Value = -50.06

$ws.Range("D17").Value = 27.1
$ws.Range("E17").Value = 80.2
$ws.Range("G17").Value = 15.14
$ws.Range("H17").Value = -47.45

$ws.Range("D18").Value = 24.01
$ws.Range("E18").Value = 71.27
$ws.Range("G18").Value = -18.920000000000002
$ws.Range("H18").Value = 71.27

$ws.Range("D19").Value = 36.32
$ws.Range("E19").Value = 53.69
$ws.Range("G19").Value = -24.19
$ws.Range("H19").Value = -36.549999999999997

# Restore the active selection to D12 and set normal-zoom view without a frozen
# top-left cell, matching the author's saved view state.
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("D12").Select() | Out-Null
